$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": insert a new client row for asesor
# "LINDAO ZUÑIGA BRYAN JOSE" (new client "IMPORTELECTRIC S.A.S") right above
# the existing "INTERNEGOCIOS DE HIERRO S.A." row (row 218), which pushes
# every following row down by one.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows("218:218").Insert()

$ws1.Range("A218").Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws1.Range("B218").Value = "IMPORTELECTRIC S.A.S"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(218, $col).Value = 0
}

# Explicit figure corrections for existing rows (all above the inserted row,
# so their row numbers are unaffected by the insert above).
$ws1.Range("M70").Value = 1258.95
$ws1.Range("P70").Value = -13.77
$ws1.Range("M99").Value = 168.38
$ws1.Range("M110").Value = 1253.61
$ws1.Range("M170").Value = -185.33
$ws1.Range("M196").Value = 15614.03

# The trailing "n de 351" counter row shifted from row 353 to row 354; update
# the denominator to reflect the new total of 352 client rows.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(354, $col)
    $txt = $cell.Text
    $cell.Value = $txt.Replace("351", "352")
}

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": mirror the same new-client insertion (row 222, same
# position relative to "INTERNEGOCIOS DE HIERRO S.A.").
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows("222:222").Insert()

$ws2.Range("A222").Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws2.Range("B222").Value = "IMPORTELECTRIC S.A.S"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(222, $col).Value = 0
}

# Explicit figure corrections for existing rows (all above the inserted row).
$ws2.Range("F70").Value = 3941.56
$ws2.Range("F99").Value = 1003.6
$ws2.Range("F110").Value = 1253.61
$ws2.Range("F174").Value = -1489.37
$ws2.Range("F200").Value = 23858.61

# Grand-total row shifted from row 357 to row 358; only the "noviembre"
# column total changes (follows from the F-column corrections above).
$ws2.Range("F358").Value = 139918.07

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": no rows inserted here, just refreshed
# VENTA / POR CUMPLIR / CUMPLIMIENTO figures for the affected asesores.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D20").Value = -13.77
$ws3.Range("E20").Value = 363.77
$ws3.Range("F20").Value = -0.03934285714285714

$ws3.Range("D24").Value = 10989.62
$ws3.Range("E24").Value = 39317.38
$ws3.Range("F24").Value = 0.2184511101834735

$ws3.Range("D36").Value = 13499.49
$ws3.Range("E36").Value = 51444.51
$ws3.Range("F36").Value = 0.2078635439763488

$ws3.Range("D48").Value = 6656.88
$ws3.Range("E48").Value = 37761.12
$ws3.Range("F48").Value = 0.1498689720383628

$ws3.Range("D60").Value = 21058.27
$ws3.Range("E60").Value = 26982.73
$ws3.Range("F60").Value = 0.4383395433067588

$ws3.Range("D77").Value = 139783.61
$ws3.Range("E77").Value = 277465.0697415455
$ws3.Range("F77").Value = 0.3350127077372314
